# New crime data collected - weekly CompStat update (84th Precinct)
# Updates the report header (volume/week number + date range) and the
# Week-to-Date / 28-Day / Year-to-Date / 2-Year crime figures table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header text: bump the report "Number" 15 -> 16 and roll the covered
#    week forward by 7 days (4/10/2023-4/16/2023 -> 4/17/2023-4/23/2023).
# ---------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value = "Volume 30   Number  16"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# ---------------------------------------------------------------------
# 2. A few cells flip between a numeric value and the "N/A" placeholder
#    text used throughout this sheet ("0" = shared string 20, "***.*" =
#    shared string 21). Plain .Value assignment of a numeric-looking
#    string collapses to a real number, so first clone the style+content
#    of an already-correct placeholder/numeric cell from an untouched
#    donor row (row 26), then overwrite with the real target value.
# ---------------------------------------------------------------------

# -> becomes "N/A" placeholder text (style 14)
$ws.Range("C26").Copy($ws.Range("D18"))
$ws.Range("E26").Copy($ws.Range("E18"))
$ws.Range("C26").Copy($ws.Range("D23"))
$ws.Range("E26").Copy($ws.Range("E23"))

# -> becomes a real number (style 15) / percent (style 16)
$ws.Range("I26").Copy($ws.Range("C22"))
$ws.Range("I26").Copy($ws.Range("C27"))
$ws.Range("I26").Copy($ws.Range("D30"))
$ws.Range("K26").Copy($ws.Range("E30"))
$ws.Range("I26").Copy($ws.Range("G30"))
$ws.Range("K26").Copy($ws.Range("H30"))

$ws.Range("C22").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0

# ---------------------------------------------------------------------
# 3. Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 27.272727272727
$ws.Range("I16").Value = 33
$ws.Range("J16").Value = 41
$ws.Range("K16").Value = -19.512195121951
$ws.Range("L16").Value = -10.810810810810
$ws.Range("M16").Value = -51.470588235294
$ws.Range("N16").Value = -90.883977900552

# ---------------------------------------------------------------------
# 4. Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 46.153846153846
$ws.Range("I17").Value = 78
$ws.Range("J17").Value = 46
$ws.Range("K17").Value = 69.565217391304
$ws.Range("L17").Value = 105.263157894737
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -35.537190082644

# ---------------------------------------------------------------------
# 5. Row 18 - Burglary (D18/E18 placeholders set above)
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 70
$ws.Range("J18").Value = 64
$ws.Range("K18").Value = 9.375
$ws.Range("L18").Value = 59.090909090909
$ws.Range("M18").Value = 141.379310344828
$ws.Range("N18").Value = -69.696969696969

# ---------------------------------------------------------------------
# 6. Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 70
$ws.Range("F19").Value = 68
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = 94.285714285714
$ws.Range("I19").Value = 209
$ws.Range("J19").Value = 178
$ws.Range("K19").Value = 17.415730337078
$ws.Range("L19").Value = 44.137931034482
$ws.Range("M19").Value = 54.814814814814
$ws.Range("N19").Value = -32.797427652733

# ---------------------------------------------------------------------
# 7. Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 6.25
$ws.Range("L20").Value = 142.857142857143
$ws.Range("M20").Value = 6.25
$ws.Range("N20").Value = -92.129629629629

# ---------------------------------------------------------------------
# 8. Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 26
$ws.Range("E21").Value = 62.5
$ws.Range("F21").Value = 122
$ws.Range("G21").Value = 73
$ws.Range("H21").Value = 67.123287671232
$ws.Range("I21").Value = 411
$ws.Range("J21").Value = 349
$ws.Range("K21").Value = 17.765042979942
$ws.Range("L21").Value = 49.454545454545
$ws.Range("M21").Value = 42.708333333333
$ws.Range("N21").Value = -67.067307692307

# ---------------------------------------------------------------------
# 9. Row 22 - Transit (C22 set above)
# ---------------------------------------------------------------------
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 15
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = 15.384615384615
$ws.Range("L22").Value = -6.25
$ws.Range("M22").Value = -31.818181818181

# ---------------------------------------------------------------------
# 10. Row 23 - Housing (D23/E23 placeholders set above)
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 6
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 15
$ws.Range("K23").Value = 15.384615384615
$ws.Range("L23").Value = 25
$ws.Range("M23").Value = 150

# ---------------------------------------------------------------------
# 11. Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 52
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 73.333333333333
$ws.Range("F24").Value = 165
$ws.Range("G24").Value = 147
$ws.Range("H24").Value = 12.244897959183
$ws.Range("I24").Value = 606
$ws.Range("J24").Value = 486
$ws.Range("K24").Value = 24.691358024691
$ws.Range("L24").Value = 45.323741007194
$ws.Range("M24").Value = 47.087378640776

# ---------------------------------------------------------------------
# 12. Row 25 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -44.444444444444
$ws.Range("F25").Value = 28
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = -12.5
$ws.Range("I25").Value = 117
$ws.Range("J25").Value = 103
$ws.Range("K25").Value = 13.592233009708
$ws.Range("L25").Value = 77.272727272727
$ws.Range("M25").Value = 0

# ---------------------------------------------------------------------
# 13. Row 27 - Other Sex Crimes (C27 set above)
# ---------------------------------------------------------------------
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 14
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = -36.363636363636
$ws.Range("L27").Value = -6.666666666666

# ---------------------------------------------------------------------
# 14. Row 30 - Hate Crimes (D30/E30/G30/H30 set above)
# ---------------------------------------------------------------------
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = 166.666666666667
